$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.474.51"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.823.84"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D5").Value = "'312.47"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4232"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'0.3625"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "'0.07191"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'0.8600"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "'20.62"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.897.82"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").Value = "'5.393"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "'6.463"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'0.06932"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'80.09"
$ws.Range("D18").Value = "'0.000008902"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'15.32"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "27.140.92"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "'5.126"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("D23").Value = "'10.90"
$ws.Range("E23").Value = "  +5.30%  "
$ws.Range("D24").Value = "1.976.97"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'1.980"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "'155.01"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "'5.154"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "'113.92"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").Value = "'1.790"
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("D31").Value = "'0.08836"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "'2.975"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").Value = "'0.7434"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").Value = "'4.533"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").Value = "'1.115"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'1.087"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("D38").Value = "'0.05267"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'2.772"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "'0.5047"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "'0.1641"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'6.445"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "'8.302"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "'10.34"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "'105.77"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").Value = "'0.4653"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D50").Value = "'1.612"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "'63.50"
$ws.Range("E51").Value = "  -1.19%  "
